$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (I1:K1), mirroring existing header strings
$ws.Range("I1").Value = "D"
$ws.Range("J1").Value = "D"
$ws.Range("K1").Value = "A"

# New data points for the extra columns
$ws.Range("I4").Value = 1
$ws.Range("J8").Value = 1

# Move the active selection to J7, matching the refreshed input layout
$ws.Range("J7").Select()
